$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# SSH public key text re-used for the moved KEY column (row 3)
$sshKey = "---- BEGIN SSH2 PUBLIC KEY ----`nComment: example`nAAAAB3NzaC1kc3MAAACBAMXXH+SzAIPRN38GehSARboF873Ic5utBjMcXx1IuFNTvvBi2j`nfCyjCBqb66sgS8pdzUl+YyH4sMDp67Q9RKI9po3ePtV03rldPQjtqgmXt2B2eQ6SKXDO3g`n+nN9LLEVXp9MpS7g9VnyDWUQCAxu+Khp+uZDzwSy7IVxRm/HHU2dAAAAFQCl1FWQ7bTyoY`n7RtEvB6rhqGyY/8QAAAIEAxYgBAFfVKvSC3AZkwWuB4hPLlBeKhL4Yt87vblimHWlaOSFU`nllKnCGmdc7R2NL3JZFP210yjapZY25YTKpkO8pdavazVqbzBd1EEtZ93umDqWua2yqPOc8`n6MoZJbk7OTJjZRlpd1XZwSI3XgyxaDtf+tCh14ikG13k4A1iKd3/MAAACBALHyHX29XFe3`nVseZeG+CiYMfc3qXbMQgpWdZeopg/1Z3qw46Kx4iiNgtZcB7BdoYdIhDvTu+xkffbG22h9`nYQnxyM9Kz/cqjKdKHp2VBX/IJU4vEkIPF+kdFPToLvJc+qkIvd1kDqUUW+e6dD6PkpNDdh`ngOn/vcgro4IwufBActyG`n---- END SSH2 PUBLIC KEY ----"

# --- Column layout change -------------------------------------------------
# The KEY and GROUPS columns swap places (KEY: H -> E, GROUPS: E -> H),
# and a brand-new COMPLIANCE_GROUPS column is appended in I.

# Header row
$ws.Range("E1").Value = "KEY"
$ws.Range("H1").Value = "GROUPS"
$ws.Range("I1").Value = "COMPLIANCE_GROUPS"

# Row 2 (server01 / SSH with password)
$ws.Range("E2").ClearContents()
$ws.Range("H2").Value = "production, test"
$ws.Range("I2").Value = "Anssi"

# Row 3 (server02 / SSH with key) - the long SSH key moves from H to E
$ws.Range("E3").Value = $sshKey
$ws.Range("E3").WrapText = $true
$ws.Range("H3").Clear()
$ws.Range("I3").Value = "Anssi"

# Row 4 (127.0.0.1 / WinRM)
$ws.Range("E4").ClearContents()
$ws.Range("H4").Value = "preproduction"
$ws.Range("I4").Value = "Anssi"

# Emphasise the relocated GROUPS column with its own explicit font
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").Font.Name = "Calibri"
$ws.Range("H1").Font.Size = 11
$ws.Range("H2").Font.Bold = $true
$ws.Range("H2").Font.Name = "Calibri"
$ws.Range("H2").Font.Size = 11
$ws.Range("H4").Font.Bold = $true
$ws.Range("H4").Font.Name = "Calibri"
$ws.Range("H4").Font.Size = 11

# --- Row heights ------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 13.8
$ws.Rows.Item(3).RowHeight = 163.5

# --- Selection (matches the author's last cursor position) -----------
$ws.Range("I3").Select()
